$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.286.10'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '2.923.10'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '369.99'
$ws.Range("E5").Value = '  +6.23%  '
$ws.Range("D6").Value = '104.32'
$ws.Range("E6").Value = '  -1.46%  '
$ws.Range("E7").Value = '  -1.47%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '0.588'
$ws.Range("E9").Value = '  -2.79%  '
$ws.Range("D10").Value = '36.81'
$ws.Range("E10").Value = '  -1.38%  '
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("D12").Value = '0.0836'
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("D13").Value = '18.36'
$ws.Range("E13").Value = '  -2.38%  '
$ws.Range("D14").Value = '3.386.85'
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("D15").Value = '7.40'
$ws.Range("E15").Value = '  -1.68%  '
$ws.Range("D16").Value = '2.928.75'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '0.941'
$ws.Range("E17").Value = '  -1.72%  '
$ws.Range("D18").Value = '51.237.68'
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").Value = '3.25'
$ws.Range("E19").Value = '  -4.55%  '
$ws.Range("D20").Value = '7.23'
$ws.Range("E20").Value = '  -1.01%  '
$ws.Range("D21").Value = '12.96'
$ws.Range("E21").Value = '  -2.69%  '
$ws.Range("D22").Value = '0.0₃0943'
$ws.Range("E22").Value = '  -1.23%  '
$ws.Range("D23").Value = '68.37'
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("D24").Value = '260.13'
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("D26").Value = '4.34'
$ws.Range("E26").Value = '  +3.76%  '
$ws.Range("E27").Value = '  +2.35%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = '25.80'
$ws.Range("E29").Value = '  -1.72%  '
$ws.Range("D30").Value = '7.06'
$ws.Range("E30").Value = '  -6.33%  '
$ws.Range("D31").Value = '0.103'
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("D32").Value = '6.27'
$ws.Range("E32").Value = '  +3.93%  '
$ws.Range("D33").Value = '9.92'
$ws.Range("E33").Value = '  -2.21%  '
$ws.Range("E34").Value = '  -1.01%  '
$ws.Range("D35").Value = '34.87'
$ws.Range("E35").Value = '  -1.14%  '
$ws.Range("D36").Value = '50.97'
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("E37").Value = '  +0.47%  '
$ws.Range("D38").Value = '0.0423'
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("D39").Value = '3.03'
$ws.Range("E39").Value = '  -1.66%  '
$ws.Range("D40").Value = '2.66'
$ws.Range("E40").Value = '  +1.75%  '
$ws.Range("E41").Value = '  -1.89%  '
$ws.Range("E42").Value = '  -4.27%  '
$ws.Range("E43").Value = '  -1.53%  '
$ws.Range("D44").Value = '22.25'
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("D45").Value = '119.21'
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("D46").Value = '2.07'
$ws.Range("E46").Value = '  -2.91%  '
$ws.Range("D47").Value = '2.024.47'
$ws.Range("E47").Value = '  -2.71%  '
$ws.Range("E48").Value = '  +2.18%  '
$ws.Range("E49").Value = '  -3.03%  '
$ws.Range("D50").Value = '0.244'
$ws.Range("E50").Value = '  +3.71%  '
$ws.Range("D51").Value = '3.216.92'
$ws.Range("E51").Value = '  +0.65%  '
